$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Widen the table / second column.
#    9207 -> 9297 dxa total table width (twips); 2160 -> 2250 dxa for the
#    second column (gridCol + every cell's tcW in that column).
#    1 point = 20 dxa, so 2250/20 = 112.5 pt and 9297/20 = 464.85 pt.
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Columns.Item(2).Width = 112.5
$t.PreferredWidth = 464.85

# ---------------------------------------------------------------------
# 2. "Fronto-central" -> "Left fronto-central" (second row, second
#    column). Insert "Left " before the word, then lower-case the
#    leading "F" to "f", keeping the existing spell-check proof markers
#    wrapping the word intact.
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Fronto-central", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($r.Start, $r.Start)
    $insertPoint.InsertBefore("Left ")

    $r2 = $d.Content
    $found2 = $r2.Find.Execute("Fronto-central", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $rF = $d.Range($r2.Start, $r2.Start + 1)
        $rF.Text = "f"
    }
}

# ---------------------------------------------------------------------
# 3. "Right parietal" -> "Right centro-parietal" (last row, second
#    column), and relocate the "_GoBack" bookmark so it now sits between
#    "centro-" and "parietal" (it previously trailed the "98" value in
#    the last column).
# ---------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("Right parietal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rpStart = $r3.Start
    $insertPoint2 = $d.Range($rpStart + 6, $rpStart + 6)
    $insertPoint2.InsertBefore("centro-")

    $bmPos = $rpStart + 6 + 7
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}
